# HighLevelSequenceDiagrams.pptx edit
# Commit: "rename: DG's image and pptx into product"
#
# The sequence diagram renamed its AddressBook-domain sample names to a
# generic "SchedulePlanner"/"Task" naming:
#   deletePerson(p)                    -> deleteTask(t)
#   AddressBookChangedEvent            -> SchedulePlannerChangedEvent
#   handleAddresssBookChangedEvent     -> handleSchedulePlannerChangedEvent
# Several label textboxes were also nudged/resized (they are left-anchored
# on arrows, so their width/position was adjusted to keep them roughly
# centered over their arrows after the text length changed).

function Replace-SubText {
    param($TextRange, [string]$OldSub, [string]$NewSub)

    $full = $TextRange.Text
    $idx = $full.IndexOf($OldSub)
    if ($idx -lt 0) {
        throw "Substring '$OldSub' not found in '$full'"
    }
    $TextRange.Characters($idx + 1, $OldSub.Length).Text = $NewSub
}

$EMU_PER_POINT = 12700.0

function Set-ShapeRectEmu {
    param($Shape, $X, $Y, $Cx, $Cy)

    $Shape.Left = $X / $EMU_PER_POINT
    $Shape.Top = $Y / $EMU_PER_POINT
    $Shape.Width = $Cx / $EMU_PER_POINT
    $Shape.Height = $Cy / $EMU_PER_POINT
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "deletePerson(p)" -> "deleteTask(t)"
$sh = $s.Shapes.Item("TextBox 28")
$tr = $sh.TextFrame.TextRange
Replace-SubText $tr "deletePerson" "deleteTask"
Replace-SubText $tr "(p)" "(t)"

# "post(AddressBookChangedEvent)" textbox next to :Model -> :EventsCenter arrow
$sh = $s.Shapes.Item("TextBox 32")
Replace-SubText $sh.TextFrame.TextRange "AddressBookChangedEvent" "SchedulePlannerChangedEvent"
Set-ShapeRectEmu $sh 5943991 1697160 2688967 215444

# "post(AddressBookChangedEvent)" textbox next to :UI -> :EventsCenter arrow
$sh = $s.Shapes.Item("TextBox 61")
Replace-SubText $sh.TextFrame.TextRange "AddressBookChangedEvent" "SchedulePlannerChangedEvent"
Set-ShapeRectEmu $sh 1717072 4797284 2716635 215444

# "handleAddresssBookChangedEvent()" textbox, :EventsCenter -> :Storage arrow (accent6)
$sh = $s.Shapes.Item("TextBox 73")
Replace-SubText $sh.TextFrame.TextRange "handleAddresssBookChangedEvent" "handleSchedulePlannerChangedEvent"
Set-ShapeRectEmu $sh 4781217 5065911 2914983 215444

# "handleAddresssBookChangedEvent()" textbox, :EventsCenter -> :UI arrow (green)
$sh = $s.Shapes.Item("TextBox 49")
Replace-SubText $sh.TextFrame.TextRange "handleAddresssBookChangedEvent" "handleSchedulePlannerChangedEvent"
Set-ShapeRectEmu $sh 1346833 5395369 2960068 215444
